$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 48037
$ws.Cells.Item(33, 9).Value = 63794.832
$ws.Cells.Item(33, 11).Value = 63794.832
$ws.Cells.Item(33, 13).Value = -63565.832
$ws.Cells.Item(64, 8).Value = 4449.9
$ws.Cells.Item(64, 9).Value = 3750
$ws.Cells.Item(64, 10).Value = 4624.875
$ws.Cells.Item(64, 11).Value = 3750
$ws.Cells.Item(64, 12).Value = 4624.875
$ws.Cells.Item(64, 13).Value = -3502
$ws.Cells.Item(64, 14).Value = -5120.875
$ws.Cells.Item(67, 8).Value = 4449.9
$ws.Cells.Item(67, 9).Value = 3750
$ws.Cells.Item(67, 10).Value = 4624.875
$ws.Cells.Item(67, 11).Value = 3750
$ws.Cells.Item(67, 12).Value = 4624.875
$ws.Cells.Item(67, 13).Value = -2892
$ws.Cells.Item(67, 14).Value = -6340.875
$ws.Cells.Item(106, 8).Value = 3734.3076
$ws.Cells.Item(106, 9).Value = 1703.6818
$ws.Cells.Item(106, 11).Value = 1703.6818
$ws.Cells.Item(106, 13).Value = -1072.6818
$ws.Cells.Item(116, 8).Value = 30761588
$ws.Cells.Item(116, 9).Value = 27895756
$ws.Cells.Item(116, 10).Value = 33340836
$ws.Cells.Item(116, 11).Value = 27895756
$ws.Cells.Item(116, 12).Value = 33340836
$ws.Cells.Item(116, 13).Value = -27892314
$ws.Cells.Item(116, 14).Value = -33347720
$ws.Cells.Item(125, 8).Value = 1174
$ws.Cells.Item(125, 9).Value = 965.55554
$ws.Cells.Item(125, 11).Value = 8689.99986
$ws.Cells.Item(125, 13).Value = -6229.99986
$ws.Cells.Item(127, 8).Value = 2059.1667
$ws.Cells.Item(127, 9).Value = 1888.75
$ws.Cells.Item(127, 11).Value = 5666.25
$ws.Cells.Item(127, 13).Value = -706.25
$ws.Cells.Item(132, 8).Value = 2328.0312
$ws.Cells.Item(132, 9).Value = 2024.1184
$ws.Cells.Item(132, 10).Value = 3482.9
$ws.Cells.Item(132, 11).Value = 6072.3552
$ws.Cells.Item(132, 12).Value = 10448.7
$ws.Cells.Item(132, 13).Value = -3542.3552
$ws.Cells.Item(132, 14).Value = -15508.7
$ws.Cells.Item(138, 8).Value = 3101.1428
$ws.Cells.Item(138, 9).Value = 2685
$ws.Cells.Item(138, 11).Value = 8055
$ws.Cells.Item(138, 13).Value = -2915
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 351.46
$ws.Cells.Item(32, 9).Value = 325.45264
$ws.Cells.Item(32, 10).Value = 845.6
$ws.Cells.Item(32, 11).Value = 325.45264
$ws.Cells.Item(32, 12).Value = 845.6
$ws.Cells.Item(32, 13).Value = -38.45263999999997
$ws.Cells.Item(32, 14).Value = -1419.6
$ws.Cells.Item(61, 8).Value = 2280.0688
$ws.Cells.Item(61, 9).Value = 1386.8096
$ws.Cells.Item(61, 11).Value = 1386.8096
$ws.Cells.Item(61, 13).Value = -1174.8096
$ws.Cells.Item(74, 8).Value = 2267.3809
$ws.Cells.Item(74, 9).Value = 2130.8125
$ws.Cells.Item(74, 11).Value = 2130.8125
$ws.Cells.Item(74, 13).Value = -1256.8125
$ws.Cells.Item(77, 8).Value = 2267.3809
$ws.Cells.Item(77, 9).Value = 2130.8125
$ws.Cells.Item(77, 11).Value = 10654.0625
$ws.Cells.Item(77, 13).Value = -6286.0625
$ws.Cells.Item(122, 8).Value = 2953.5
$ws.Cells.Item(122, 9).Value = 1651.1818
$ws.Cells.Item(122, 11).Value = 4953.5454
$ws.Cells.Item(122, 13).Value = -2503.5454
$ws.Cells.Item(132, 8).Value = 691174.6
$ws.Cells.Item(132, 9).Value = 418050.62
$ws.Cells.Item(132, 11).Value = 1254151.86
$ws.Cells.Item(132, 13).Value = -1251621.86
$ws.Cells.Item(136, 8).Value = 2280.0688
$ws.Cells.Item(136, 9).Value = 1386.8096
$ws.Cells.Item(136, 11).Value = 4160.4288
$ws.Cells.Item(136, 13).Value = -1610.4288
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 1547.931
$ws.Cells.Item(134, 9).Value = 1428.2142
$ws.Cells.Item(134, 11).Value = 4284.642599999999
$ws.Cells.Item(134, 13).Value = -1749.642599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 6000
$ws.Cells.Item(23, 10).Value = 6000
$ws.Cells.Item(23, 12).Value = 6000
$ws.Cells.Item(23, 14).Value = -6480
$ws.Cells.Item(27, 8).Value = 6000
$ws.Cells.Item(27, 10).Value = 6000
$ws.Cells.Item(27, 12).Value = 6000
$ws.Cells.Item(27, 14).Value = -6384
$ws.Cells.Item(31, 8).Value = 2889.1538
$ws.Cells.Item(31, 9).Value = 1802
$ws.Cells.Item(31, 11).Value = 1802
$ws.Cells.Item(31, 13).Value = -1507
$ws.Cells.Item(34, 8).Value = 2889.1538
$ws.Cells.Item(34, 9).Value = 1802
$ws.Cells.Item(34, 11).Value = 1802
$ws.Cells.Item(34, 13).Value = -1600
$ws.Cells.Item(86, 8).Value = 59712.535
$ws.Cells.Item(86, 9).Value = 78997
$ws.Cells.Item(86, 10).Value = 6680.25
$ws.Cells.Item(86, 11).Value = 78997
$ws.Cells.Item(86, 12).Value = 6680.25
$ws.Cells.Item(86, 13).Value = -77874
$ws.Cells.Item(86, 14).Value = -8926.25
$ws.Cells.Item(89, 8).Value = 59712.535
$ws.Cells.Item(89, 9).Value = 78997
$ws.Cells.Item(89, 10).Value = 6680.25
$ws.Cells.Item(89, 11).Value = 394985
$ws.Cells.Item(89, 12).Value = 33401.25
$ws.Cells.Item(89, 13).Value = -389369
$ws.Cells.Item(89, 14).Value = -44633.25
$ws.Cells.Item(94, 8).Value = 2770.9473
$ws.Cells.Item(94, 9).Value = 3154.7778
$ws.Cells.Item(94, 10).Value = 2425.5
$ws.Cells.Item(94, 11).Value = 3154.7778
$ws.Cells.Item(94, 12).Value = 2425.5
$ws.Cells.Item(94, 13).Value = -2703.7778
$ws.Cells.Item(94, 14).Value = -3327.5
$ws.Cells.Item(134, 8).Value = 2537.3667
$ws.Cells.Item(134, 9).Value = 1807.1111
$ws.Cells.Item(134, 10).Value = 3632.75
$ws.Cells.Item(134, 11).Value = 5421.3333
$ws.Cells.Item(134, 12).Value = 10898.25
$ws.Cells.Item(134, 13).Value = -2886.3333
$ws.Cells.Item(134, 14).Value = -15968.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 737.8570999999999
$ws.Cells.Item(5, 9).Value = 327.2857
$ws.Cells.Item(5, 10).Value = 1148.4286
$ws.Cells.Item(5, 11).Value = 981.8571000000001
$ws.Cells.Item(5, 12).Value = 3445.2858
$ws.Cells.Item(5, 13).Value = -869.8571000000001
$ws.Cells.Item(5, 14).Value = -3669.2858
$ws.Cells.Item(122, 8).Value = 255.9
$ws.Cells.Item(122, 10).Value = 195
$ws.Cells.Item(122, 12).Value = 1755
$ws.Cells.Item(122, 14).Value = -6655
$ws.Cells.Item(131, 8).Value = 14875.8125
$ws.Cells.Item(131, 9).Value = 842.1429000000001
$ws.Cells.Item(131, 10).Value = 25790.889
$ws.Cells.Item(131, 11).Value = 2526.4287
$ws.Cells.Item(131, 12).Value = 77372.667
$ws.Cells.Item(131, 13).Value = 2513.5713
$ws.Cells.Item(131, 14).Value = -87452.667
$ws.Cells.Item(135, 8).Value = 737.8570999999999
$ws.Cells.Item(135, 9).Value = 327.2857
$ws.Cells.Item(135, 10).Value = 1148.4286
$ws.Cells.Item(135, 11).Value = 2945.5713
$ws.Cells.Item(135, 12).Value = 10335.8574
$ws.Cells.Item(135, 13).Value = -410.5713000000001
$ws.Cells.Item(135, 14).Value = -15405.8574
$ws.Cells.Item(138, 8).Value = 8587261
$ws.Cells.Item(138, 9).Value = 20000610
$ws.Cells.Item(138, 11).Value = 60001830
$ws.Cells.Item(138, 13).Value = -59996690
$ws.Cells.Item(140, 8).Value = 1952.5454
$ws.Cells.Item(140, 9).Value = 1434.75
$ws.Cells.Item(140, 11).Value = 4304.25
$ws.Cells.Item(140, 13).Value = 875.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 9996.5
$ws.Cells.Item(2, 9).Value = 9995.333000000001
$ws.Cells.Item(2, 11).Value = 9995.333000000001
$ws.Cells.Item(2, 13).Value = -9882.333000000001
$ws.Cells.Item(107, 8).Value = 688.4545000000001
$ws.Cells.Item(107, 9).Value = 299.16666
$ws.Cells.Item(107, 10).Value = 1155.6
$ws.Cells.Item(107, 11).Value = 299.16666
$ws.Cells.Item(107, 12).Value = 1155.6
$ws.Cells.Item(107, 13).Value = 1620.83334
$ws.Cells.Item(107, 14).Value = -4995.6
$ws.Cells.Item(113, 8).Value = 3801.389
$ws.Cells.Item(113, 9).Value = 2701.7144
$ws.Cells.Item(113, 11).Value = 2701.7144
$ws.Cells.Item(113, 13).Value = -531.7143999999998
$ws.Cells.Item(132, 8).Value = 361042.53
$ws.Cells.Item(132, 9).Value = 479568.66
$ws.Cells.Item(132, 10).Value = 5464.143
$ws.Cells.Item(132, 11).Value = 1438705.98
$ws.Cells.Item(132, 12).Value = 16392.429
$ws.Cells.Item(132, 13).Value = -1436175.98
$ws.Cells.Item(132, 14).Value = -21452.429
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 950000
$ws.Cells.Item(4, 9).Value = 1000000
$ws.Cells.Item(4, 10).Value = 900000
$ws.Cells.Item(4, 11).Value = 1000000
$ws.Cells.Item(4, 12).Value = 900000
$ws.Cells.Item(4, 13).Value = -999887
$ws.Cells.Item(4, 14).Value = -900226
$ws.Cells.Item(22, 8).Value = 2305609.5
$ws.Cells.Item(22, 10).Value = 6454393
$ws.Cells.Item(22, 12).Value = 6454393
$ws.Cells.Item(22, 14).Value = -6454983
$ws.Cells.Item(27, 8).Value = 2305609.5
$ws.Cells.Item(27, 10).Value = 6454393
$ws.Cells.Item(27, 12).Value = 6454393
$ws.Cells.Item(27, 14).Value = -6454607
$ws.Cells.Item(28, 8).Value = 950000
$ws.Cells.Item(28, 9).Value = 1000000
$ws.Cells.Item(28, 10).Value = 900000
$ws.Cells.Item(28, 11).Value = 1000000
$ws.Cells.Item(28, 12).Value = 900000
$ws.Cells.Item(28, 13).Value = -999768
$ws.Cells.Item(28, 14).Value = -900464
$ws.Cells.Item(37, 8).Value = 950000
$ws.Cells.Item(37, 9).Value = 1000000
$ws.Cells.Item(37, 10).Value = 900000
$ws.Cells.Item(37, 11).Value = 1000000
$ws.Cells.Item(37, 12).Value = 900000
$ws.Cells.Item(37, 13).Value = -999893
$ws.Cells.Item(37, 14).Value = -900214
$ws.Cells.Item(40, 8).Value = 3712.7334
$ws.Cells.Item(40, 9).Value = 3069.2
$ws.Cells.Item(40, 11).Value = 3069.2
$ws.Cells.Item(40, 13).Value = -2933.2
$ws.Cells.Item(46, 8).Value = 2638.1155
$ws.Cells.Item(46, 9).Value = 1966.4445
$ws.Cells.Item(46, 11).Value = 1966.4445
$ws.Cells.Item(46, 13).Value = -1778.4445
$ws.Cells.Item(82, 8).Value = 1535.5714
$ws.Cells.Item(82, 9).Value = 1033
$ws.Cells.Item(82, 10).Value = 1912.5
$ws.Cells.Item(82, 11).Value = 1033
$ws.Cells.Item(82, 12).Value = 1912.5
$ws.Cells.Item(82, 13).Value = -672
$ws.Cells.Item(82, 14).Value = -2634.5
$ws.Cells.Item(85, 8).Value = 1535.5714
$ws.Cells.Item(85, 9).Value = 1033
$ws.Cells.Item(85, 10).Value = 1912.5
$ws.Cells.Item(85, 11).Value = 1033
$ws.Cells.Item(85, 12).Value = 1912.5
$ws.Cells.Item(85, 13).Value = 215
$ws.Cells.Item(85, 14).Value = -4408.5
$ws.Cells.Item(93, 8).Value = 1045.0952
$ws.Cells.Item(93, 9).Value = 943.1667
$ws.Cells.Item(93, 11).Value = 943.1667
$ws.Cells.Item(93, 13).Value = 304.8333
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 1500
$ws.Cells.Item(28, 10).Value = 1500
$ws.Cells.Item(28, 12).Value = 1500
$ws.Cells.Item(28, 14).Value = -2196
$ws.Cells.Item(100, 8).Value = 283.5
$ws.Cells.Item(100, 9).Value = 280.7143
$ws.Cells.Item(100, 11).Value = 561.4286
$ws.Cells.Item(100, 13).Value = -20.42859999999996
$ws.Cells.Item(132, 8).Value = 230180.94
$ws.Cells.Item(132, 9).Value = 265933.25
$ws.Cells.Item(132, 11).Value = 797799.75
$ws.Cells.Item(132, 13).Value = -795269.75
